$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated price/volume figures and two coin-row swaps as per the
# automated "Updated cryptos list" GitHub Actions run.

$ws.Range("D2").Value = '41.715.47'
$ws.Range("E2").Value = '  +0.16%  '
$ws.Range("D3").Value = '2.209.77'
$ws.Range("E3").Value = '  -1.95%  '
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '229.96'
$ws.Range("E5").Value = '  -2.16%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.615'
$ws.Range("E6").Value = '  -4.11%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '60.32'
$ws.Range("E7").Value = '  -5.08%  '
$ws.Range("E8").Value = '  -0.14%  '
$ws.Range("E9").Value = '  -1.79%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '57.23'
$ws.Range("E10").Value = '  -4.26%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0890'
$ws.Range("E11").Value = '  -0.59%  '
$ws.Range("E12").Value = '  -2.00%  '
$ws.Range("D13").Value = '2.536.82'
$ws.Range("E13").Value = '  -2.18%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '15.42'
$ws.Range("E14").Value = '  -3.97%  '
$ws.Range("E15").Value = '  -2.61%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.66'
$ws.Range("E17").Value = '  -3.61%  '
$ws.Range("D18").Value = '2.213.94'
$ws.Range("E18").Value = '  -2.50%  '
$ws.Range("D19").Value = '41.697.30'
$ws.Range("E19").Value = '  +0.33%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '72.15'
$ws.Range("E20").Value = '  -3.32%  '
$ws.Range("D21").Value = '0.0₃0901'
$ws.Range("E21").Value = '  -3.15%  '
$ws.Range("E22").Value = '  -1.52%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '242.50'
$ws.Range("E23").Value = '  -3.61%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.998'
$ws.Range("E24").Value = '  -0.28%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.48'
$ws.Range("E25").Value = '  +1.93%  '
$ws.Range("E26").Value = '  -3.18%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.68'
$ws.Range("E27").Value = '  -1.67%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '169.21'
$ws.Range("E28").Value = '  -1.23%  '
$ws.Range("E29").Value = '  -5.13%  '
$ws.Range("E30").Value = '  +0.40%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '19.75'
$ws.Range("E31").Value = '  -3.44%  '
$ws.Range("E32").Value = '  -7.76%  '
$ws.Range("E33").Value = '  -3.67%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.03'
$ws.Range("E34").Value = '  -1.53%  '
$ws.Range("E35").Value = '  -3.23%  '
$ws.Range("E36").Value = '  +1.93%  '
$ws.Range("E37").Value = '  -3.50%  '
$ws.Range("E38").Value = '  -7.56%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.55'
$ws.Range("E39").Value = '  -7.89%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.000240'
$ws.Range("E40").Value = '  -9.39%  '
$ws.Range("E41").Value = '  -0.26%  '
$ws.Range("E42").Value = '  -0.47%  '
$ws.Range("E43").Value = '  -1.24%  '
$ws.Range("B44").Value = 'FTXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.43'
$ws.Range("E44").Value = '  -12.63%  '
$ws.Range("B45").Value = 'Cronos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0954'
$ws.Range("E45").Value = '  -3.46%  '
$ws.Range("B46").Value = 'TrustWalletToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.20'
$ws.Range("E46").Value = '  -2.98%  '
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '97.32'
$ws.Range("E47").Value = '  -4.75%  '
$ws.Range("D48").Value = '1.466.53'
$ws.Range("E48").Value = '  -2.79%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '16.39'
$ws.Range("E49").Value = '  -6.90%  '
$ws.Range("E50").Value = '  -1.47%  '
$ws.Range("E51").Value = '  -4.95%  '

Write-Host "Applied cryptos list update"
